$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new worksheet right after "总计" for the new "2022-Q4"
#    quarter, pushing all the other quarter sheets one slot to the right.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

# headers
$q4.Cells.Item(1,2).Value = "基金代码"
$q4.Cells.Item(1,3).Value = "基金名称"
$q4.Cells.Item(1,4).Value = "基金规模"
$q4.Cells.Item(1,5).Value = "股票总仓位"
$q4.Cells.Item(1,6).Value = "仓位占比"
$q4.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4.Cells.Item(1,8).Value = "仓位排名"

# row 2 - 南方宝裕混合A
$q4.Cells.Item(2,1).Value = 0
$q4.Cells.Item(2,2).Value = "012945"
$q4.Cells.Item(2,3).Value = "南方宝裕混合A"
$q4.Cells.Item(2,4).Value = "11.23"
$q4.Cells.Item(2,5).Value = "22.43"
$q4.Cells.Item(2,6).Value = "0.59"
$q4.Cells.Item(2,7).Value = "0.0663"
$q4.Cells.Item(2,8).Value = 10

# row 3 - 南方宝裕混合C
$q4.Cells.Item(3,1).Value = 1
$q4.Cells.Item(3,2).Value = "012946"
$q4.Cells.Item(3,3).Value = "南方宝裕混合C"
$q4.Cells.Item(3,4).Value = "0.23"
$q4.Cells.Item(3,5).Value = "22.43"
$q4.Cells.Item(3,6).Value = "0.59"
$q4.Cells.Item(3,7).Value = "0.0014"
$q4.Cells.Item(3,8).Value = 10

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert a new row 2 for 2022-Q4 and
#    shift the existing quarters down by one row.
# ---------------------------------------------------------------------------
$total.Rows.Item(2).Insert()

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 2
$total.Cells.Item(2,4).Value = 0.07

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2022-Q3"
$total.Cells.Item(3,3).Value = 4
$total.Cells.Item(3,4).Value = 0.63

$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = "2022-Q2"
$total.Cells.Item(4,3).Value = 6
$total.Cells.Item(4,4).Value = 0.32

$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(5,2).Value = "2022-Q1"
$total.Cells.Item(5,3).Value = 13
$total.Cells.Item(5,4).Value = 0.96

$total.Cells.Item(6,1).Value = 4
$total.Cells.Item(6,2).Value = "2021-Q4"
$total.Cells.Item(6,3).Value = 12
$total.Cells.Item(6,4).Value = 1.29

$total.Cells.Item(7,1).Value = 5
$total.Cells.Item(7,2).Value = "2021-Q3"
$total.Cells.Item(7,3).Value = 3
$total.Cells.Item(7,4).Value = 0.75

$total.Cells.Item(8,1).Value = 6
$total.Cells.Item(8,2).Value = "2021-Q1"
$total.Cells.Item(8,3).Value = 4
$total.Cells.Item(8,4).Value = 0.1

$total.Cells.Item(9,1).Value = 7
$total.Cells.Item(9,2).Value = "2020-Q4"
$total.Cells.Item(9,3).Value = 1
$total.Cells.Item(9,4).Value = 0.1
